$d = $word.ActiveDocument

# Locate the paragraph holding "Ver no Jupiter Salvar em pdf Salvar em docx".
# Together with the blank paragraph just before it, the blank paragraph just
# after it, and the following page-break paragraph, these four paragraphs
# are being removed (the trailing blank + page-break paragraphs further
# down stay untouched).
$rng = $d.Content
$found = $rng.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$matchIndex = $rng.Paragraphs.Item(1).Index

$startPara = $d.Paragraphs.Item($matchIndex - 1)
$endPara = $d.Paragraphs.Item($matchIndex + 2)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
